$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 8754
$ws.Range("B2").Value = "Bryan Mendes"
$ws.Range("C2").Value = "Juridico"
$ws.Range("D2").Value = "Problemas pessoais"
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 45102
$ws.Range("G2").Value = 7738.22

# Row 3
$ws.Range("A3").Value = 68229
$ws.Range("B3").Value = "Arthur Gabriel Viana"
$ws.Range("C3").Value = "Atendimento ao Cliente"
$ws.Range("D3").Value = "Problemas pessoais"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 45080
$ws.Range("G3").Value = 8289.42

# Row 4
$ws.Range("A4").Value = 53883
$ws.Range("B4").Value = "Pedro Nascimento"
$ws.Range("C4").Value = "Vendas"
$ws.Range("D4").Value = "Doenca"
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 45094
$ws.Range("G4").Value = 2886.63

# Row 5
$ws.Range("A5").Value = 81817
$ws.Range("B5").Value = "Sr. Joaquim Moura"
$ws.Range("C5").Value = "P&D"
$ws.Range("D5").Value = "Doenca"
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 45088
$ws.Range("G5").Value = 8046.97

# Row 6
$ws.Range("A6").Value = 2099
$ws.Range("B6").Value = "Sr. André Oliveira"
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 45102
$ws.Range("G6").Value = 9205.940000000001

# Row 7
$ws.Range("A7").Value = 43707
$ws.Range("B7").Value = "Bento Guerra"
$ws.Range("C7").Value = "Operacoes"
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 45095
$ws.Range("G7").Value = 7798.02

# Row 8
$ws.Range("A8").Value = 51070
$ws.Range("B8").Value = "Helena Cardoso"
$ws.Range("C8").Value = "Engenharia"
$ws.Range("D8").Value = "Viagem de negocios"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 45091
$ws.Range("G8").Value = 2734.44

# Row 9
$ws.Range("A9").Value = 94006
$ws.Range("B9").Value = "Alexandre Costela"
$ws.Range("C9").Value = "Engenharia"
$ws.Range("D9").Value = "Outros"
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 45087
$ws.Range("G9").Value = 2094.35

# Row 10
$ws.Range("A10").Value = 90732
$ws.Range("B10").Value = "Murilo Siqueira"
$ws.Range("C10").Value = "Atendimento ao Cliente"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 45088
$ws.Range("G10").Value = 2393.3

# Row 11
$ws.Range("A11").Value = 68124
$ws.Range("B11").Value = "Vicente Macedo"
$ws.Range("C11").Value = "Marketing"
$ws.Range("D11").Value = "Doenca"
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 45101
$ws.Range("G11").Value = 5494.27
